$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Jason Roy"
$ws.Range("B2").Value = 11
$ws.Range("C2").Value = 4
$ws.Range("E2").Value = " Haris Rauf"
$ws.Range("J2").Value = "Mohammad Rizwan"
$ws.Range("K2").Value = 47
$ws.Range("L2").Value = 15
$ws.Range("M2").Value = "Bowled"
$ws.Range("N2").Value = " Chris Woakes"
# Row 3
$ws.Range("A3").Value = "Jos Buttler"
$ws.Range("B3").Value = 4
$ws.Range("D3").Value = "Caught"
$ws.Range("E3").Value = " Hasan Ali"
$ws.Range("J3").Value = "Babar Azam(C)"
$ws.Range("K3").Value = 47
$ws.Range("L3").Value = 16
$ws.Range("M3").Value = "LBW"
$ws.Range("N3").Value = " Chris Woakes"
# Row 4
$ws.Range("A4").Value = "Dawid Malan"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "Bowled"
$ws.Range("E4").Value = " Haris Rauf"
$ws.Range("J4").Value = "Fakhar Zaman"
$ws.Range("K4").Value = 10
$ws.Range("L4").Value = 3
$ws.Range("M4").Value = "LBW"
$ws.Range("N4").Value = " Chris Woakes"
# Row 5
$ws.Range("A5").Value = "Jonny Bairstow"
$ws.Range("B5").Value = 18
$ws.Range("C5").Value = 6
$ws.Range("E5").Value = " Shaheen Afridi"
$ws.Range("J5").Value = "Mohammad Hafeez"
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 1
$ws.Range("N5").Value = " Chris Woakes"
# Row 6
$ws.Range("A6").Value = "Eoin Morgan(C)"
$ws.Range("B6").Value = 31
$ws.Range("C6").Value = 13
$ws.Range("D6").Value = "Bowled"
$ws.Range("E6").Value = " Hasan Ali"
$ws.Range("J6").Value = "Shoaib Malik"
$ws.Range("K6").Value = 11
$ws.Range("L6").Value = 7
$ws.Range("N6").Value = " Mark Wood"
# Row 7
$ws.Range("A7").Value = "Moeen Ali"
$ws.Range("B7").Value = 27
$ws.Range("C7").Value = 14
$ws.Range("D7").Value = "LBW"
$ws.Range("E7").Value = " Haris Rauf"
$ws.Range("J7").Value = "Asif Ali"
$ws.Range("K7").Value = 35
$ws.Range("L7").Value = 15
$ws.Range("M7").Value = "* NOT OUT"
$ws.Range("N7").Value = " "
# Row 8
$ws.Range("A8").Value = "Liam Livingstone"
$ws.Range("B8").Value = 18
$ws.Range("C8").Value = 6
$ws.Range("D8").Value = "NOT OUT"
$ws.Range("E8").Value = " "
$ws.Range("J8").Value = "Shadab Khan"
$ws.Range("K8").Value = 14
$ws.Range("L8").Value = 7
$ws.Range("M8").Value = "LBW"
$ws.Range("N8").Value = " Mark Wood"
# Row 9
$ws.Range("A9").Value = "Chris Woakes"
$ws.Range("B9").Value = 27
$ws.Range("C9").Value = 8
$ws.Range("E9").Value = " Imad Wasim"
$ws.Range("J9").Value = "Imad Wasim"
$ws.Range("K9").Value = 1
$ws.Range("M9").Value = "NOT OUT"
# Row 10
$ws.Range("A10").Value = "Chris Jordan"
$ws.Range("B10").Value = 7
$ws.Range("C10").Value = 4
$ws.Range("E10").Value = " Shadab Khan"
$ws.Range("J10").Value = "Hasan Ali"
# Row 11
$ws.Range("A11").Value = "Adil Rashid"
$ws.Range("B11").Value = 19
$ws.Range("C11").Value = 8
$ws.Range("E11").Value = " Shaheen Afridi"
$ws.Range("J11").Value = "Shaheen Afridi"
# Row 12
$ws.Range("A12").Value = "Mark Wood"
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = "Caught"
$ws.Range("E12").Value = " Shaheen Afridi"
$ws.Range("J12").Value = "Haris Rauf"
# Row 16
$ws.Range("A16").Value = 162
$ws.Range("C16").Value = "'11.2"
$ws.Range("D16").Value = 68
$ws.Range("J16").Value = 165
$ws.Range("L16").Value = "'11.0"
$ws.Range("M16").Value = 66
# Row 21
$ws.Range("A21").Value = "Hasan Ali"
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = 13.5
$ws.Range("J21").Value = "Adil Rashid"
$ws.Range("L21").Value = 42
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 21
# Row 22
$ws.Range("A22").Value = "Imad Wasim"
$ws.Range("B22").Value = "'2.0"
$ws.Range("C22").Value = 37
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 18.5
$ws.Range("J22").Value = "Chris Jordan"
$ws.Range("L22").Value = 24
$ws.Range("N22").Value = 12
# Row 23
$ws.Range("A23").Value = "Shadab Khan"
$ws.Range("B23").Value = "'2.0"
$ws.Range("C23").Value = 24
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 12
$ws.Range("J23").Value = "Chris Woakes"
$ws.Range("L23").Value = 23
$ws.Range("M23").Value = 4
$ws.Range("N23").Value = 11.5
# Row 24
$ws.Range("A24").Value = "Haris Rauf"
$ws.Range("C24").Value = 40
$ws.Range("D24").Value = 3
$ws.Range("E24").Value = 13.33
$ws.Range("J24").Value = "Liam Livingstone"
$ws.Range("K24").Value = "'2.0"
$ws.Range("L24").Value = 26
$ws.Range("M24").Value = 0
$ws.Range("N24").Value = 13
# Row 25
$ws.Range("A25").Value = "Shaheen Afridi"
$ws.Range("B25").Value = "'2.2"
$ws.Range("C25").Value = 34
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 15.45
$ws.Range("J25").Value = "Mark Wood"
$ws.Range("K25").Value = "'3.0"
$ws.Range("L25").Value = 50
$ws.Range("N25").Value = 16.67

Write-Host "Applied all changes"